# Adds the "issues found on Aug 16, 2022" rows to three sheets:
#   termsWithMulitpleLabels  (sheet1) -> new section at rows 256-268
#   LabelsUsedMultipleTerms  (sheet2) -> new section at rows 131-133
#   termWithDifferentParent  (sheet3) -> new section at rows 280-281
#
# Cell formats are cloned from existing same-style cells via
# Copy()/PasteSpecial(xlPasteFormats) so the workbook's style table
# (cellXfs/fonts) is reused rather than growing with near-duplicate entries.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("termsWithMulitpleLabels")
$ws2 = $wb.Worksheets.Item("LabelsUsedMultipleTerms")
$ws3 = $wb.Worksheets.Item("termWithDifferentParent")

$xlPasteFormats = -4122

# ======================================================================
# Sheet "termsWithMulitpleLabels": new section header + table at A256:D268
# ======================================================================

$ws1.Cells.Item(256,1).Value = "Checked on Aug 16th, 2022, identified following issue:"
$ws1.Cells.Item(254,1).Copy()
$ws1.Cells.Item(256,1).PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# header row
$ws1.Cells.Item(257,1).Value = "sid"
$ws1.Cells.Item(257,2).Value = "label_count"
$ws1.Cells.Item(257,3).Value = "labels"
$ws1.Cells.Item(257,4).Value = "datasets"

# data rows
$ws1.Cells.Item(258,1).Value = "EUPATH_0053214"
$ws1.Cells.Item(258,2).Value = 2
$ws1.Cells.Item(258,3).Value = "Sprayed against mosquitoes in past 12 months | Sprayed in past 12 months"
$ws1.Cells.Item(258,4).Value = "gates_namibia | icemr_sw_pacific"

$ws1.Cells.Item(259,1).Value = "EUPATH_0053590"
$ws1.Cells.Item(259,2).Value = 2
$ws1.Cells.Item(259,3).Value = "Was there sexual behavior change since take PrEP | Sexual behavior change since take PrEP"
$ws1.Cells.Item(259,4).Value = "gates_jilinde_prospective_cohort | gates_jilinde_retrospective_survey"

$ws1.Cells.Item(260,1).Value = "EUPATH_0053594"
$ws1.Cells.Item(260,2).Value = 2
$ws1.Cells.Item(260,3).Value = "Reason to stop taking PrEP | Reason you stopped taking PrEP"
$ws1.Cells.Item(260,4).Value = "gates_jilinde_awareness_survey | gates_jilinde_prospective_cohort | gates_jilinde_retrospective_survey | gates_jilinde_demand_creation_evaluation_questionnaire"

$ws1.Cells.Item(261,1).Value = "EUPATH_0054079"
$ws1.Cells.Item(261,2).Value = 2
$ws1.Cells.Item(261,3).Value = "What did you do with extra tablets after stopping PrEP | Extra tablets after stop taking PrEP"
$ws1.Cells.Item(261,4).Value = "gates_jilinde_prospective_cohort | gates_jilinde_retrospective_survey"

$ws1.Cells.Item(262,1).Value = "EUPATH_0054080"
$ws1.Cells.Item(262,2).Value = 2
$ws1.Cells.Item(262,3).Value = "Other use of extra tablets after stopping PrEP specified | Extra tablets after stop taking PrEP specified"
$ws1.Cells.Item(262,4).Value = "gates_jilinde_prospective_cohort | gates_jilinde_retrospective_survey"

$ws1.Cells.Item(263,1).Value = "EUPATH_0054162"
$ws1.Cells.Item(263,2).Value = 2
$ws1.Cells.Item(263,3).Value = "Live together with spouse or partner | Live with partner"
$ws1.Cells.Item(263,4).Value = "gates_jilinde_awareness_survey | gates_jilinde_prospective_cohort | gates_jilinde_demand_creation_evaluation_questionnaire"

$ws1.Cells.Item(264,1).Value = "EUPATH_0054169"
$ws1.Cells.Item(264,2).Value = 2
$ws1.Cells.Item(264,3).Value = "Other main reason for agreeing to offer of PrEP specified | Main reason for agreeing to offer of PrEP"
$ws1.Cells.Item(264,4).Value = "gates_jilinde_prospective_cohort"

$ws1.Cells.Item(265,1).Value = "EUPATH_0054175"
$ws1.Cells.Item(265,2).Value = 2
$ws1.Cells.Item(265,3).Value = "Main reason for stopping PrEP after prior use | Other main reason for stopping PrEP after prior use"
$ws1.Cells.Item(265,4).Value = "gates_jilinde_prospective_cohort"

$ws1.Cells.Item(266,1).Value = "EUPATH_0054248"
$ws1.Cells.Item(266,2).Value = 2
$ws1.Cells.Item(266,3).Value = "Other reason of reconsidering taking PrEP specified | Reason of reconsidering taking PrEP"
$ws1.Cells.Item(266,4).Value = "gates_jilinde_prospective_cohort"

$ws1.Cells.Item(267,1).Value = "EUPATH_0054330"
$ws1.Cells.Item(267,2).Value = 2
$ws1.Cells.Item(267,3).Value = "Whether have child | Have children"
$ws1.Cells.Item(267,4).Value = "gates_jilinde_awareness_survey | gates_jilinde_prospective_cohort | gates_jilinde_demand_creation_evaluation_questionnaire"

$ws1.Cells.Item(268,1).Value = "EUPATH_0054344"
$ws1.Cells.Item(268,2).Value = 2
$ws1.Cells.Item(268,3).Value = "Sex work | Commercial sex and activity"
$ws1.Cells.Item(268,4).Value = "gates_jilinde_prospective_cohort | gates_jilinde_retrospective_survey | gates_jilinde_demand_creation_evaluation_questionnaire"

# ======================================================================
# Sheet "LabelsUsedMultipleTerms": new section header + row at A131:G133
# ======================================================================

$ws2.Cells.Item(131,1).Value = "Checked on Aug 16th, 2022, one issue is found"
$ws2.Cells.Item(128,1).Copy()
$ws2.Cells.Item(131,1).PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws2.Cells.Item(132,1).Value = "Enrollment date"
$ws2.Cells.Item(132,2).Value = 2
$ws2.Cells.Item(132,3).Value = "EUPATH_0000151|EUPATH_0054076"
$ws2.Cells.Item(132,4).Value = "variable"
$ws2.Cells.Item(132,5).Value = "Administrative information"
$ws2.Cells.Item(132,6).Value = "changed to EUPATH_0000151"
$ws2.Cells.Item(124,6).Copy()
$ws2.Cells.Item(132,6).PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws2.Cells.Item(132,7).Value = "Fixed"

# row 133 is an otherwise-empty row that only carries a touched F cell
# (matches the pre-existing "touched-then-cleared" placeholder pattern
# already present elsewhere in this sheet/workbook, e.g. F127/F128)
$ws2.Cells.Item(133,6).Value = "x"
$ws2.Cells.Item(133,6).Value = $null

# ======================================================================
# Sheet "termWithDifferentParent": new section header + row at A280:F281
# ======================================================================

$ws3.Cells.Item(280,1).Value = "Checked on Aug 16th, 2022, one issue is found."
$ws3.Cells.Item(278,1).Copy()
$ws3.Cells.Item(280,1).PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws3.Cells.Item(281,1).Value = "EUPATH_0053590"
$ws3.Cells.Item(281,2).Value = 2
$ws3.Cells.Item(281,3).Value = "Was there sexual behavior change since take PrEP|Sexual behavior change since take PrEP"
$ws3.Cells.Item(281,4).Value = "PrEP|Sexual behavior"
$ws3.Cells.Item(281,5).Value = "gates_jilinde_prospective_cohort | gates_jilinde_retrospective_survey"
$ws3.Cells.Item(281,6).Value = "x"
$ws3.Cells.Item(281,6).Value = $null

# ======================================================================
# Window / selection state: active sheet moves from LabelsUsedMultipleTerms
# to termsWithMulitpleLabels; each sheet keeps its own last selection.
# Non-active sheets' selections are set first; the to-be-active sheet is
# activated and selected last so its tabSelected flag "sticks".
# ======================================================================

$ws2.Range("A133").Select()
$ws3.Range("C284").Select()

$ws1.Activate()
$ws1.Range("C257").Select()
